# Add a new "tahirih" / 12345 test-data row to Sheet1, and leave Sheet1 as
# the active sheet/selection (matches the commit's updated testData.xlsx).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A6").Value = "tahirih"
$ws1.Range("B6").Value = 12345

# Select Sheet1!B6 last so it becomes the active sheet/tab (activeTab -> 0)
# and the active selection (mirrors the workbookView/sheetView changes in
# the diff: Sheet1 tabSelected=true, Sheet3 tabSelected=false).
$ws1.Activate()
$ws1.Range("B6").Select()
